$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Common (unchanged) column values for every data row in this sheet.
$colA = 3
$colB = "Femacal de La Calera"
$colC = "Coquimbo"
$colE = 5
$colF = 300000000
$colG = "Esp" + [char]0x00E1 + "rragos"
$colH = "Verde"
$colN = "$/kilo"
$colQ = 1
$colR = "Hortaliza"

function Fill-Row($r, $date, $quality, $j, $k, $l, $m, $origin, $p) {
    $ws.Cells.Item($r, 1).Value = $colA
    $ws.Cells.Item($r, 2).Value = $colB
    $ws.Cells.Item($r, 3).Value = $colC
    $ws.Cells.Item($r, 4).Value = $date
    $ws.Cells.Item($r, 5).Value = $colE
    $ws.Cells.Item($r, 6).Value = $colF
    $ws.Cells.Item($r, 7).Value = $colG
    $ws.Cells.Item($r, 8).Value = $colH
    $ws.Cells.Item($r, 9).Value = $quality
    $ws.Cells.Item($r, 10).Value = $j
    $ws.Cells.Item($r, 11).Value = $k
    $ws.Cells.Item($r, 12).Value = $l
    $ws.Cells.Item($r, 13).Value = $m
    $ws.Cells.Item($r, 14).Value = $colN
    $ws.Cells.Item($r, 15).Value = $origin
    $ws.Cells.Item($r, 16).Value = $p
    $ws.Cells.Item($r, 17).Value = $colQ
    $ws.Cells.Item($r, 18).Value = $colR
}

# Insert a new row at position 5 — a new weekly price record (Provincia de
# Linares), pushing the former rows 5..37 down to 6..38.
$ws.Rows.Item(5).Insert()
Fill-Row 5 44859 "Primera" 1580 1400 1400 1400 "Provincia de Linares" 1400

# Insert a second new row, now at position 38 (right before the former row
# 37, which has since shifted to row 38 and will end up at row 39).
$ws.Rows.Item(38).Insert()
Fill-Row 38 44858 "Primera" 2400 1300 1400 1346 "Provincia de Quillota" 1346
